# Sync attendance_reports from main repo:
# Column G ("Recorded By") values of "System, dnasr281@gmail.com" are
# re-ordered to "dnasr281@gmail.com, System" throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "System, dnasr281@gmail.com"
$newVal = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldVal) {
        $cell.Value = $newVal
    }
}
